$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 4.077281088115026
$ws.Range("E2").Value = 3.013456857613561
$ws.Range("C3").Value = 4.731933890736406
$ws.Range("E3").Value = 4.482374336309247
$ws.Range("C4").Value = 7.427427238257622
$ws.Range("E4").Value = 6.028771162940338
$ws.Range("C5").Value = 3.697244903694119
$ws.Range("E5").Value = 5.35733476656457
$ws.Range("C6").Value = -0.9010546343133807
$ws.Range("E6").Value = 1.483308540745609
$ws.Range("C7").Value = 1.278852728916302
$ws.Range("E7").Value = 1.211545622441634
$ws.Range("C8").Value = 2.40113223806655
$ws.Range("E8").Value = 1.221134982408678
$ws.Range("C9").Value = 0.4170416928886977
$ws.Range("E9").Value = 1.535040327807513
$ws.Range("C10").Value = 2.025199397970145
$ws.Range("E10").Value = 1.766027057877517
$ws.Range("C11").Value = 2.535130037318867
$ws.Range("E11").Value = 2.27519030414034
$ws.Range("C12").Value = 1.085017960020163
$ws.Range("E12").Value = 1.594198396297974
$ws.Range("C13").Value = 2.388449447315399
$ws.Range("E13").Value = 1.896063224966515
$ws.Range("C14").Value = 1.940295589655605
$ws.Range("E14").Value = 2.152498180268503
$ws.Range("C15").Value = 0.1294783794713039
$ws.Range("E15").Value = 1.357330623126884
$ws.Range("C16").Value = 0.2239492534813481
$ws.Range("E16").Value = 0.1832633470655098
$ws.Range("C17").Value = 0.669497318959178
$ws.Range("E17").Value = 0.4274836977099516
$ws.Range("C18").Value = 0.9511721486624936
$ws.Range("E18").Value = 0.886066450790457
$ws.Range("C19").Value = 1.350833417525776
$ws.Range("E19").Value = 1.249180524815863
$ws.Range("C20").Value = 3.305550968939119
$ws.Range("E20").Value = 2.444559947892744
$ws.Range("C21").Value = 2.741128804567849
$ws.Range("E21").Value = 3.207064487734335
$ws.Range("C22").Value = -5.478868953971427
$ws.Range("E22").Value = -2.015486574969738
$ws.Range("C23").Value = -0.1094048593225039
$ws.Range("E23").Value = -1.306092631642397
$ws.Range("C24").Value = 3.381937564063731
$ws.Range("E24").Value = 0.8627271536207459
$ws.Range("C25").Value = 1.637009187238481
$ws.Range("E25").Value = 2.774647569643585
$ws.Range("C26").Value = -0.03183845066089264
$ws.Range("E26").Value = 1.102201969172678
$ws.Range("C27").Value = 1.810762846774527
$ws.Range("E27").Value = 1.207964354105195
$ws.Range("C28").Value = 1.292027346513414
$ws.Range("E28").Value = 0.9897401519578963
$ws.Range("C29").Value = 1.566699735746391
$ws.Range("E29").Value = 1.595944879287448
$ws.Range("C30").Value = 1.752318341645176
$ws.Range("E30").Value = 1.832607040765044
$ws.Range("C31").Value = 2.349173111882341
$ws.Range("E31").Value = 2.133017022807637
$ws.Range("C32").Value = 0.8252516928923814
$ws.Range("E32").Value = 1.704160690624001
$ws.Range("C33").Value = -1.611564731980897
$ws.Range("E33").Value = -0.741724804865096
$ws.Range("C34").Value = -0.5453046728852495
$ws.Range("E34").Value = -1.767889269204159
$ws.Range("C35").Value = 1.726013280798222
$ws.Range("E35").Value = 0.05925426804285205
$ws.Range("C36").Value = -0.06641493770841445
$ws.Range("E36").Value = 0.7991555368092929
$ws.Range("C37").Value = -0.1397017661237232
$ws.Range("E37").Value = 0.385156833908451
$ws.Range("C38").Value = 0.06464796496492564
$ws.Range("E38").Value = -0.06351196001971315
